{"js": "// Update the date line and every three-digit-by-one-digit multiplication\n// answer in the table to the new values from the target revision.\nconst replacements = [\n  [\"2025-08-29 Friday\", \"2025-08-30 Saturday\"],\n  [\"113\u00d72=226\", \"268\u00d75=1340\"],\n  [\"875\u00d79=7875\", \"794\u00d79=7146\"],\n  [\"399\u00d74=1596\", \"621\u00d77=4347\"],\n  [\"459\u00d73=1377\", \"250\u00d76=1500\"],\n  [\"432\u00d73=1296\", \"679\u00d72=1358\"],\n  [\"146\u00d72=292\", \"624\u00d78=4992\"],\n  [\"972\u00d75=4860\", \"529\u00d79=4761\"],\n  [\"941\u00d73=2823\", \"466\u00d75=2330\"],\n  [\"372\u00d73=1116\", \"336\u00d78=2688\"],\n  [\"683\u00d73=2049\", \"536\u00d74=2144\"],\n  [\"820\u00d76=4920\", \"291\u00d73=873\"],\n  [\"703\u00d77=4921\", \"613\u00d72=1226\"],\n  [\"997\u00d78=7976\", \"257\u00d75=1285\"],\n  [\"408\u00d72=816\", \"593\u00d74=2372\"],\n  [\"958\u00d74=3832\", \"597\u00d74=2388\"],\n  [\"423\u00d79=3807\", \"275\u00d75=1375\"],\n  [\"209\u00d78=1672\", \"112\u00d76=672\"],\n  [\"422\u00d76=2532\", \"856\u00d79=7704\"],\n  [\"991\u00d75=4955\", \"456\u00d76=2736\"],\n  [\"394\u00d73=1182\", \"439\u00d75=2195\"],\n  [\"986\u00d79=8874\", \"922\u00d72=1844\"],\n  [\"348\u00d73=1044\", \"354\u00d74=1416\"],\n  [\"403\u00d77=2821\", \"736\u00d76=4416\"],\n  [\"300\u00d72=600\", \"582\u00d72=1164\"],\n  [\"732\u00d77=5124\", \"150\u00d74=600\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and every three-digit-by-one-digit multiplication\n# answer in the table to the new values from the target revision.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-08-29 Friday\", \"2025-08-30 Saturday\"),\n    @(\"113\u00d72=226\", \"268\u00d75=1340\"),\n    @(\"875\u00d79=7875\", \"794\u00d79=7146\"),\n    @(\"399\u00d74=1596\", \"621\u00d77=4347\"),\n    @(\"459\u00d73=1377\", \"250\u00d76=1500\"),\n    @(\"432\u00d73=1296\", \"679\u00d72=1358\"),\n    @(\"146\u00d72=292\", \"624\u00d78=4992\"),\n    @(\"972\u00d75=4860\", \"529\u00d79=4761\"),\n    @(\"941\u00d73=2823\", \"466\u00d75=2330\"),\n    @(\"372\u00d73=1116\", \"336\u00d78=2688\"),\n    @(\"683\u00d73=2049\", \"536\u00d74=2144\"),\n    @(\"820\u00d76=4920\", \"291\u00d73=873\"),\n    @(\"703\u00d77=4921\", \"613\u00d72=1226\"),\n    @(\"997\u00d78=7976\", \"257\u00d75=1285\"),\n    @(\"408\u00d72=816\", \"593\u00d74=2372\"),\n    @(\"958\u00d74=3832\", \"597\u00d74=2388\"),\n    @(\"423\u00d79=3807\", \"275\u00d75=1375\"),\n    @(\"209\u00d78=1672\", \"112\u00d76=672\"),\n    @(\"422\u00d76=2532\", \"856\u00d79=7704\"),\n    @(\"991\u00d75=4955\", \"456\u00d76=2736\"),\n    @(\"394\u00d73=1182\", \"439\u00d75=2195\"),\n    @(\"986\u00d79=8874\", \"922\u00d72=1844\"),\n    @(\"348\u00d73=1044\", \"354\u00d74=1416\"),\n    @(\"403\u00d77=2821\", \"736\u00d76=4416\"),\n    @(\"300\u00d72=600\", \"582\u00d72=1164\"),\n    @(\"732\u00d77=5124\", \"150\u00d74=600\")\n)\n\nforeach ($pair in $replacements) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
